$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1481.5714
$ws.Range("I28").Value = 1750.2
$ws.Range("J28").Value = 810
$ws.Range("K28").Value = 1750.2
$ws.Range("L28").Value = 810
$ws.Range("M28").Value = -1265.2
$ws.Range("N28").Value = -1780
$ws.Range("H113").Value = 3347.6924
$ws.Range("I113").Value = 2683.75
$ws.Range("J113").Value = 4410
$ws.Range("K113").Value = 2683.75
$ws.Range("L113").Value = 4410
$ws.Range("M113").Value = 570.25
$ws.Range("N113").Value = -10918
$ws.Range("H132").Value = 6582707.5
$ws.Range("I132").Value = 3787.6365
$ws.Range("J132").Value = 50003580
$ws.Range("K132").Value = 11362.9095
$ws.Range("L132").Value = 150010740
$ws.Range("M132").Value = -8832.9095
$ws.Range("N132").Value = -150015800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1421716.2
$ws.Range("I45").Value = 1749478.5
$ws.Range("J45").Value = 1413.3334
$ws.Range("K45").Value = 1749478.5
$ws.Range("L45").Value = 1413.3334
$ws.Range("M45").Value = -1749101.5
$ws.Range("N45").Value = -2167.3334
$ws.Range("H61").Value = 2433.814
$ws.Range("I61").Value = 2324.4614
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 2324.4614
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -2112.4614
$ws.Range("N61").Value = -3924
$ws.Range("H74").Value = 4640.6665
$ws.Range("I74").Value = 4640.6665
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4640.6665
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 4640.6665
$ws.Range("I77").Value = 4640.6665
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 23203.3325
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -5244
$ws.Range("H122").Value = 10392.6875
$ws.Range("I122").Value = 11413.071
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 34239.213
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -31789.213
$ws.Range("N122").Value = -14650
$ws.Range("H127").Value = 59091.285
$ws.Range("J127").Value = 58988.332
$ws.Range("L127").Value = 58988.332
$ws.Range("N127").Value = -68908.33199999999
$ws.Range("H129").Value = 48999.75
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 48999.75
$ws.Range("K129").Value = 0
$ws.Range("L129").ClearContents()
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -58999.75
$ws.Range("H132").Value = 1753.0364
$ws.Range("I132").Value = 1242.6595
$ws.Range("J132").Value = 4751.5
$ws.Range("K132").Value = 3727.9785
$ws.Range("L132").Value = 14254.5
$ws.Range("M132").Value = -1197.9785
$ws.Range("N132").Value = -19314.5
$ws.Range("H136").Value = 2433.814
$ws.Range("I136").Value = 2324.4614
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 6973.3842
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -4423.3842
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2018.3684
$ws.Range("I20").Value = 1998.3334
$ws.Range("J20").Value = 2093.5
$ws.Range("K20").Value = 1998.3334
$ws.Range("L20").Value = 2093.5
$ws.Range("M20").Value = -1751.3334
$ws.Range("N20").Value = -2587.5
$ws.Range("H53").Value = 47503.332
$ws.Range("J53").Value = 47503.332
$ws.Range("L53").Value = 47503.332
$ws.Range("N53").Value = -48651.332
$ws.Range("H99").Value = 1222.2222
$ws.Range("I99").Value = 1125
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1125
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 373
$ws.Range("N99").Value = -4996
$ws.Range("H115").Value = 37508.668
$ws.Range("J115").Value = 37508.668
$ws.Range("L115").Value = 37508.668
$ws.Range("N115").Value = -40642.668
$ws.Range("H129").Value = 49834.75
$ws.Range("J129").Value = 49834.75
$ws.Range("L129").Value = 49834.75
$ws.Range("N129").Value = -59834.75
$ws.Range("H134").Value = 3492.3442
$ws.Range("I134").Value = 2313.6858
$ws.Range("J134").Value = 5079
$ws.Range("K134").Value = 6941.057400000001
$ws.Range("L134").Value = 15237
$ws.Range("M134").Value = -4406.057400000001
$ws.Range("N134").Value = -20307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4822.2563
$ws.Range("I31").Value = 4135.6665
$ws.Range("J31").Value = 5127.407
$ws.Range("K31").Value = 4135.6665
$ws.Range("L31").Value = 5127.407
$ws.Range("M31").Value = -3840.6665
$ws.Range("N31").Value = -5717.407
$ws.Range("H34").Value = 4822.2563
$ws.Range("I34").Value = 4135.6665
$ws.Range("J34").Value = 5127.407
$ws.Range("K34").Value = 4135.6665
$ws.Range("L34").Value = 5127.407
$ws.Range("M34").Value = -3933.6665
$ws.Range("N34").Value = -5531.407
$ws.Range("H58").Value = 2711.8823
$ws.Range("I58").Value = 966.8333
$ws.Range("J58").Value = 6900
$ws.Range("K58").Value = 966.8333
$ws.Range("L58").Value = 6900
$ws.Range("M58").Value = -763.8333
$ws.Range("N58").Value = -7306
$ws.Range("H94").Value = 3403.9656
$ws.Range("I94").Value = 2355.9092
$ws.Range("J94").Value = 4044.4443
$ws.Range("K94").Value = 2355.9092
$ws.Range("L94").Value = 4044.4443
$ws.Range("M94").Value = -1904.9092
$ws.Range("N94").Value = -4946.4443
$ws.Range("H99").Value = 2402.45
$ws.Range("I99").Value = 1766.6666
$ws.Range("K99").Value = 1766.6666
$ws.Range("M99").Value = -268.6666
$ws.Range("H122").Value = 3655.8333
$ws.Range("I122").Value = 3587
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10761
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8311
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 2402.45
$ws.Range("I126").Value = 1766.6666
$ws.Range("K126").Value = 5299.9998
$ws.Range("M126").Value = -2829.9998
$ws.Range("H127").Value = 53731.25
$ws.Range("J127").Value = 53731.25
$ws.Range("L127").Value = 53731.25
$ws.Range("N127").Value = -63651.25
$ws.Range("H136").Value = 2711.8823
$ws.Range("I136").Value = 966.8333
$ws.Range("J136").Value = 6900
$ws.Range("K136").Value = 2900.4999
$ws.Range("L136").Value = 20700
$ws.Range("M136").Value = -350.4998999999998
$ws.Range("N136").Value = -25800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 7259.5
$ws.Range("I97").Value = 10654.9
$ws.Range("J97").Value = 1600.5
$ws.Range("K97").Value = 31964.7
$ws.Range("L97").Value = 4801.5
$ws.Range("M97").Value = -31468.7
$ws.Range("N97").Value = -5793.5
$ws.Range("H132").Value = 1350.5555
$ws.Range("I132").Value = 753.1539
$ws.Range("J132").Value = 1905.2858
$ws.Range("K132").Value = 6778.3851
$ws.Range("L132").Value = 17147.5722
$ws.Range("M132").Value = -4248.3851
$ws.Range("N132").Value = -22207.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3540.2307
$ws.Range("I122").Value = 3241.1667
$ws.Range("J122").Value = 3796.5715
$ws.Range("K122").Value = 9723.500100000001
$ws.Range("L122").Value = 11389.7145
$ws.Range("M122").Value = -7273.500100000001
$ws.Range("N122").Value = -16289.7145
$ws.Range("H126").Value = 3430.1724
$ws.Range("I126").Value = 2988.889
$ws.Range("J126").Value = 3628.75
$ws.Range("K126").Value = 8966.667000000001
$ws.Range("L126").Value = 10886.25
$ws.Range("M126").Value = -6496.667000000001
$ws.Range("N126").Value = -15826.25
$ws.Range("H132").Value = 6920.2
$ws.Range("I132").Value = 7816.316
$ws.Range("J132").Value = 4082.5
$ws.Range("K132").Value = 23448.948
$ws.Range("L132").Value = 12247.5
$ws.Range("M132").Value = -20918.948
$ws.Range("N132").Value = -17307.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10612.294
$ws.Range("I7").Value = 8501.571
$ws.Range("K7").Value = 8501.571
$ws.Range("M7").Value = -8389.571
$ws.Range("H100").Value = 2692.5715
$ws.Range("I100").Value = 2949
$ws.Range("J100").Value = 2590
$ws.Range("K100").Value = 2949
$ws.Range("L100").Value = 2590
$ws.Range("M100").Value = -2408
$ws.Range("N100").Value = -3672
$ws.Range("H122").Value = 14883.111
$ws.Range("I122").Value = 50000
$ws.Range("J122").Value = 10493.5
$ws.Range("K122").Value = 150000
$ws.Range("L122").Value = 31480.5
$ws.Range("M122").Value = -147550
$ws.Range("N122").Value = -36380.5
$ws.Range("H126").Value = 10612.294
$ws.Range("I126").Value = 8501.571
$ws.Range("K126").Value = 25504.713
$ws.Range("M126").Value = -23034.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 37285.6
$ws.Range("J46").Value = 37285.6
$ws.Range("L46").Value = 37285.6
$ws.Range("N46").Value = -37747.6
$ws.Range("H122").Value = 5284.3335
$ws.Range("I122").Value = 5441.6
$ws.Range("J122").Value = 4498
$ws.Range("K122").Value = 16324.8
$ws.Range("L122").Value = 13494
$ws.Range("M122").Value = -13874.8
$ws.Range("N122").Value = -18394
$ws.Range("H132").Value = 1582.7213
$ws.Range("I132").Value = 1348.8572
$ws.Range("J132").Value = 2537.6667
$ws.Range("K132").Value = 4046.5716
$ws.Range("L132").Value = 7613.000100000001
$ws.Range("M132").Value = -1516.5716
$ws.Range("N132").Value = -12673.0001
$ws.Range("H134").Value = 37285.6
$ws.Range("J134").Value = 37285.6
$ws.Range("L134").Value = 111856.8
$ws.Range("N134").Value = -116926.8
$ws.Range("H136").Value = 1236.986
$ws.Range("I136").Value = 594.1429000000001
$ws.Range("J136").Value = 6299.375
$ws.Range("K136").Value = 1782.4287
$ws.Range("L136").Value = 18898.125
$ws.Range("M136").Value = 767.5712999999998
$ws.Range("N136").Value = -23998.125
